$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell A2 value from "vm-test-1" to "vm-openssl-test"
$ws.Range("A2").Value = "vm-openssl-test"

# Update selection to A2 (was C2)
$ws.Range("A2").Select()
